# Update cryptocurrency price (D) and 1h volume change % (E) columns
# Commit: 'Updated cryptos list on Thu Aug 24 22:36:39 UTC 2023 with GitHub Actions'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format first so numeric-looking values
# (e.g. "218.88") are stored as text like the rest of the column, matching
# the original inlineStr cells instead of being auto-coerced to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.249.17"
$ws.Range("D3").Value = "1.662.60"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "218.88"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "0.5222"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.2671"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.06330"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "21.10"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").Value = "0.07723"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "4.431"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "1.645.63"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "1.889.85"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "0.5477"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "0.0₅8244"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "64.95"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "26.279.74"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "4.661"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "194.15"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "10.15"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "6.082"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "138.71"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").Value = "0.1241"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "7.235"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").Value = "16.13"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "1.409"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "0.05977"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").Value = "1.283"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "3.626"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "3.313"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "1.632"
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("D35").Value = "0.9798"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").Value = "2.416"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "2.786"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "0.5898"
$ws.Range("D39").Value = "0.01593"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").Value = "5.947"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "0.8618"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "1.030.35"
$ws.Range("E43").Value = "  -4.12%  "
$ws.Range("D44").Value = "99.76"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "1.803.19"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  +4.53%  "
$ws.Range("D47").Value = "57.29"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "8.113"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "0.05184"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "1.475"
$ws.Range("E51").Value = "  -0.28%  "

# Restore the default (Normal) style so no stray number-format style lingers
# on the cells -- only their text content has changed, same as the diff.
$dataRange.Style = "Normal"
